$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.330464124679565
$ws.Range("B1").Value = 1.701269865036011
$ws.Range("C1").Value = 3.357109069824219
$ws.Range("D1").Value = 3.72789478302002
$ws.Range("E1").Value = 1.264052510261536
